# Seed investigation report: insert a new experiment row (row 6) with a
# "Special setups" note column, per commit "feat: Save the excel report".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the current row 6 ("Adam"/seed=1 row), so the
#    new data point becomes row 6 and everything below shifts down by one.
# ---------------------------------------------------------------------
$ws.Rows("6:6").Insert(-4121)

# ---------------------------------------------------------------------
# 2. Fill the new row 6 with the new data point.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "CodeGPTPy"
$ws.Range("B6").Value = 128
$ws.Range("C6").Value = "Early Stopping"
$ws.Range("D6").Value = "Adam"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.70799999999999996
$ws.Range("G6").Value = 0.68799999999999994
$ws.Range("H6").Formula = "=F6-G6"

# Formatting: thin border all around + centered, matching the other
# data rows (A:E use the plain "General" number format, F:H use 0.000).
$dataRow = $ws.Range("A6:H6")
$dataRow.HorizontalAlignment = -4108   # xlCenter
$dataRow.VerticalAlignment = -4108     # xlCenter
foreach ($edge in 7, 8, 9, 10) {
    $b = $dataRow.Borders.Item($edge)
    $b.LineStyle = 1        # xlContinuous
    $b.Weight = 2           # xlThin
    $b.ColorIndex = 0
}
$ws.Range("A6:E6").NumberFormat = "General"
$ws.Range("F6:H6").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# 3. Add the new "Special setups" column (I) with a note for the new row.
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "Special setups"
$ws.Range("I6").Value = "minl2=1e-4,minl1=1e-4"
$ws.Columns("I").ColumnWidth = 20.6640625

# ---------------------------------------------------------------------
# 4. Restore selection to G6 (where the new value was entered).
# ---------------------------------------------------------------------
$ws.Range("G6").Select()
